$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels in row 8 (D8, E8)
$ws.Range("D8").Value = "breedte"
$ws.Range("E8").Value = "hoogte"

# Update D9 value (breedte value), E9 keeps its formula and recalculates
$ws.Range("D9").Value = 900
$ws.Range("E9").Formula = "=D9/5*7"
